$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (WireFrame/Kanban shift right to E/F)
$ws.Columns.Item(4).Insert()

# New header for the inserted column
$ws.Range("D1").Value = "Valeur"

# Row 2 (Page se connecter 1) - all statuses become "Fait"
$ws.Range("B2:F2").Value = "Fait"

# Remaining data rows - all statuses become "Fait" except column E which becomes lowercase "fait"
$rows = @(3, 4, 5, 7, 8, 10, 11, 13, 14, 16, 17)
foreach ($r in $rows) {
    $ws.Range("B$r").Value = "Fait"
    $ws.Range("C$r").Value = "Fait"
    $ws.Range("D$r").Value = "Fait"
    $ws.Range("E$r").Value = "fait"
    $ws.Range("F$r").Value = "Fait"
}

# Restore the previously-selected cell selection similar to the authored change
$ws.Range("D20").Select()

$wb.Save()
